$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "ワークフローファイル"
$ws.Range("B1").Value = "内部パス"
$ws.Range("C1").Value = "対象"
$ws.Range("D1").Value = "事象"
$ws.Range("E1").Value = "メッセージ"
